# Weekly fruit/vegetable price update for "Ajo" (Garlic) sheet.
# Inserts two new rows (122:123) above the existing data block, shifting the
# previously-existing rows 122-139 down to 124-141, and populates the two new
# rows with the latest weekly report values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 122:123, pushing everything from the old row 122
# downward (old row 122 -> new row 124, ..., old row 139 -> new row 141).
$ws.Rows("122:123").Insert()

# New row 122: Ajo / Chino / Primera, $/caja 10 kilos
$ws.Range("A122").Value = 9
$ws.Range("B122").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C122").Value = "Metropolitana"
$ws.Range("D122").Value = 44522
$ws.Range("E122").Value = 13
$ws.Range("F122").Value = 100112003
$ws.Range("G122").Value = "Ajo"
$ws.Range("H122").Value = "Chino"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 520
$ws.Range("K122").Value = 17000
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = 17500
$ws.Range("N122").Value = "`$/caja 10 kilos"
$ws.Range("O122").Value = "China"
$ws.Range("P122").Value = 1750
$ws.Range("Q122").Value = 10
$ws.Range("R122").Value = "Hortaliza"

# New row 123: Ajo / Chino / Primera, $/malla 10 kilos
$ws.Range("A123").Value = 9
$ws.Range("B123").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C123").Value = "Metropolitana"
$ws.Range("D123").Value = 44522
$ws.Range("E123").Value = 13
$ws.Range("F123").Value = 100112003
$ws.Range("G123").Value = "Ajo"
$ws.Range("H123").Value = "Chino"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 340
$ws.Range("K123").Value = 16000
$ws.Range("L123").Value = 17000
$ws.Range("M123").Value = 16500
$ws.Range("N123").Value = "`$/malla 10 kilos"
$ws.Range("O123").Value = "China"
$ws.Range("P123").Value = 1650
$ws.Range("Q123").Value = 10
$ws.Range("R123").Value = "Hortaliza"
